# Included Jira IDs for each test case data.
#
# Adds a new "RefID" column at the front of the sheet (holding a Jira ID
# for each existing test-case row) and a new second test-case row that
# exercises the Document Register / review-document attachment fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the left edge - this shifts the existing
# To/CC/Subject/... columns from A:O to B:P and carries their styles,
# values and column widths along with them.
$ws.Columns("A:A").Insert()

# Existing (row 2) test case now gets its Jira reference id.
$ws.Range("A2").Value = "LATFLD-37"

# New test case row 3 - same base transmittal fields as row 2, plus
# AttachDocuments/AttachSupportDocuments + the Document Register entries
# showing the review document being picked up.
$ws.Range("B3").Value = "AutoTestAdmin"
$ws.Range("C3").Value = "AutoTestUser"
$ws.Range("D3").Value = "New Transmittal from Automation"
$ws.Range("E3").Value = "UnTick"
$ws.Range("F3").Value = "Correspondence"
$ws.Range("G3").Value = "Issued for Information"
$ws.Range("H3").Value = "Document Register"
$ws.Range("I3").Value = "Test 1 ta.docx"
$ws.Range("J3").Value = "Document Register"
$ws.Range("K3").Value = "Test 1 ta.docx"
$ws.Range("M3").Value = "Message for New transmittal"
$ws.Range("A3").Value = "LATFLD-3"

# Header for the new RefID column - match the bold header style used by
# the rest of row 1.
$ws.Range("A1").Value = "RefID"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Re-fit the new/changed columns to their content, same as Excel does
# automatically when a column's data changes.
$ws.Columns("A:A").ColumnWidth = 9.0834
$ws.Columns("H:H").ColumnWidth = 17.25
$ws.Columns("I:I").ColumnWidth = 20.5834
$ws.Columns("J:J").ColumnWidth = 23.25
$ws.Columns("K:K").ColumnWidth = 27.75

$null = $ws.Range("B5").Select()
